# Insert two new weekly price rows for "Provincia de Cardenal Caro" right
# before the current row 62, shifting the existing rows (old 62-79) down
# to rows 64-81 and extending the used range to A1:T81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("62:63").Insert()

# New row 62: "Primera" quality, $/bandeja 2 kilos unit
$ws.Cells.Item(62, 1).Value = 10
$ws.Cells.Item(62, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(62, 3).Value = "La Araucanía"
$ws.Cells.Item(62, 4).Value = 44553
$ws.Cells.Item(62, 5).Value = 9
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100101
$ws.Cells.Item(62, 8).Value = "Berries"
$ws.Cells.Item(62, 9).Value = 100101001
$ws.Cells.Item(62, 10).Value = "Arándano (blue)"
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 325
$ws.Cells.Item(62, 14).Value = 2200
$ws.Cells.Item(62, 15).Value = 2500
$ws.Cells.Item(62, 16).Value = 2315
$ws.Cells.Item(62, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(62, 18).Value = "Provincia de Cardenal Caro"
$ws.Cells.Item(62, 19).Value = 1158
$ws.Cells.Item(62, 20).Value = 2

# New row 63: "Segunda" quality, $/bandeja 2 kilos unit
$ws.Cells.Item(63, 1).Value = 10
$ws.Cells.Item(63, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value = "La Araucanía"
$ws.Cells.Item(63, 4).Value = 44553
$ws.Cells.Item(63, 5).Value = 9
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100101
$ws.Cells.Item(63, 8).Value = "Berries"
$ws.Cells.Item(63, 9).Value = 100101001
$ws.Cells.Item(63, 10).Value = "Arándano (blue)"
$ws.Cells.Item(63, 11).Value = "Sin especificar"
$ws.Cells.Item(63, 12).Value = "Segunda"
$ws.Cells.Item(63, 13).Value = 150
$ws.Cells.Item(63, 14).Value = 1800
$ws.Cells.Item(63, 15).Value = 1800
$ws.Cells.Item(63, 16).Value = 1800
$ws.Cells.Item(63, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(63, 18).Value = "Provincia de Cardenal Caro"
$ws.Cells.Item(63, 19).Value = 900
$ws.Cells.Item(63, 20).Value = 2

Write-Output "Inserted rows 62-63; dimension now extends through row 81."
